# Update "想去人数" (interested-count) figures across sheets, as output
# was regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 204
$ws1.Range("F4").Value  = 406
$ws1.Range("F5").Value  = 986
$ws1.Range("F6").Value  = 5474
$ws1.Range("F7").Value  = 486
$ws1.Range("F8").Value  = 682
$ws1.Range("F9").Value  = 946
$ws1.Range("F12").Value = 35
$ws1.Range("F14").Value = 27
$ws1.Range("F15").Value = 21
$ws1.Range("F17").Value = 1825
$ws1.Range("F18").Value = 1465
$ws1.Range("F19").Value = 894
$ws1.Range("F20").Value = 297
$ws1.Range("F22").Value = 330
$ws1.Range("F23").Value = 539
$ws1.Range("F24").Value = 146
$ws1.Range("F28").Value = 2842
$ws1.Range("F29").Value = 176
$ws1.Range("F33").Value = 34
$ws1.Range("F34").Value = 367
$ws1.Range("F39").Value = 287
$ws1.Range("F40").Value = 699
$ws1.Range("F43").Value = 55

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 183
$ws2.Range("F6").Value = 128

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 204
$ws4.Range("F5").Value  = 986
$ws4.Range("F7").Value  = 5474
$ws4.Range("F8").Value  = 486
$ws4.Range("F9").Value  = 682
$ws4.Range("F11").Value = 183
$ws4.Range("F12").Value = 946
$ws4.Range("F15").Value = 128
$ws4.Range("F17").Value = 35
$ws4.Range("F19").Value = 27
$ws4.Range("F20").Value = 21
$ws4.Range("F23").Value = 1825
$ws4.Range("F24").Value = 1465
$ws4.Range("F25").Value = 895
$ws4.Range("F27").Value = 330
$ws4.Range("F29").Value = 539
$ws4.Range("F30").Value = 146
$ws4.Range("F32").Value = 2842
$ws4.Range("F33").Value = 176
$ws4.Range("F37").Value = 34
$ws4.Range("F38").Value = 367
$ws4.Range("F42").Value = 287
$ws4.Range("F43").Value = 699
$ws4.Range("F45").Value = 55
